$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("AG2").Value = 12
$ws.Range("AJ2").Value = 67
$ws.Range("AK2").Value = 401
$ws.Range("AL2").Value = 7.5
$ws.Range("AM2").Value = 7.5
$ws.Range("K2").Value = 2.38
$ws.Range("L2").Value = 2
$ws.Range("Y2").Value = 1.95
$ws.Range("Z2").Value = 1.8

# Row 3
$ws.Range("AA3").Value = 5.5
$ws.Range("AB3").Value = 9.5
$ws.Range("AG3").Value = 5.5
$ws.Range("AI3").Value = 21
$ws.Range("AL3").Value = 6.5
$ws.Range("AN3").Value = 13
$ws.Range("G3").Value = 2.4
$ws.Range("I3").Value = 3.3
$ws.Range("J3").Value = 3.4
$ws.Range("L3").Value = 4.33
$ws.Range("W3").Value = 1.67
$ws.Range("X3").Value = 2.1
$ws.Range("Y3").Value = 2.38
$ws.Range("Z3").Value = 1.53

# Row 4
$ws.Range("S4").Value = 3.6
$ws.Range("T4").Value = 1.3

# Row 6
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5
$ws.Range("O6").Value = 1.57
$ws.Range("P6").Value = 2.25

# Row 7
$ws.Range("AB7").Value = 34
$ws.Range("AD7").Value = 81
$ws.Range("AH7").Value = 7.5
$ws.Range("AN7").Value = 9
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 8

# Row 8
$ws.Range("AH8").Value = 6.5
$ws.Range("AL8").Value = 10
$ws.Range("AM8").Value = 26
$ws.Range("AN8").Value = 21
$ws.Range("AR8").Value = 2.05
$ws.Range("AS8").Value = 1.8
$ws.Range("H8").Value = 3.1
$ws.Range("I8").Value = 5.75
$ws.Range("K8").Value = 1.91
$ws.Range("L8").Value = 6
$ws.Range("M8").Value = 1.13
$ws.Range("N8").Value = 6
$ws.Range("Q8").Value = 2.7
$ws.Range("R8").Value = 1.44
$ws.Range("W8").Value = 1.62
$ws.Range("X8").Value = 2.2

# Row 9
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 8

# Row 10
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 7

# Row 11
$ws.Range("AA11").Value = 5.5
$ws.Range("AB11").Value = 13
$ws.Range("AD11").Value = 40
$ws.Range("AF11").Value = 75
$ws.Range("AG11").Value = 3.9
$ws.Range("AH11").Value = 5.2
$ws.Range("AI11").Value = 22
$ws.Range("AL11").Value = 5.8
$ws.Range("AM11").Value = 14
$ws.Range("AQ11").Value = 75
$ws.Range("G11").Value = 2.95
$ws.Range("H11").Value = 2.37
$ws.Range("I11").Value = 3.1
$ws.Range("J11").Value = 3.8
$ws.Range("K11").Value = 1.7
$ws.Range("L11").Value = 3.9
$ws.Range("N11").Value = 3.9
$ws.Range("O11").Value = 1.82
$ws.Range("P11").Value = 1.9
$ws.Range("Q11").Value = 3.4
$ws.Range("R11").Value = 1.28
$ws.Range("U11").Value = 6.3
$ws.Range("Y11").Value = 2.42
$ws.Range("Z11").Value = 1.5

# Row 12
$ws.Range("AA12").Value = 6.5
$ws.Range("AB12").Value = 15
$ws.Range("AG12").Value = 4.45
$ws.Range("AI12").Value = 19.5
$ws.Range("AL12").Value = 5.6
$ws.Range("AM12").Value = 11.25
$ws.Range("AN12").Value = 10.75
$ws.Range("G12").Value = 3.25
$ws.Range("H12").Value = 2.55
$ws.Range("I12").Value = 2.62
$ws.Range("K12").Value = 1.78
$ws.Range("L12").Value = 3.4
$ws.Range("N12").Value = 4.45
$ws.Range("O12").Value = 1.65
$ws.Range("P12").Value = 2.1
$ws.Range("Q12").Value = 2.95
$ws.Range("R12").Value = 1.35
$ws.Range("U12").Value = 5.4
$ws.Range("V12").Value = 1.12

# Row 16
$ws.Range("AA16").Value = 5
$ws.Range("AB16").Value = 6.5
$ws.Range("AD16").Value = 12
$ws.Range("G16").Value = 1.67
$ws.Range("I16").Value = 5.75
$ws.Range("J16").Value = 2.38

# Row 17
$ws.Range("AA17").Value = 7
$ws.Range("AB17").Value = 15
$ws.Range("AC17").Value = 13
$ws.Range("AD17").Value = 41
$ws.Range("AE17").Value = 34
$ws.Range("AG17").Value = 6
$ws.Range("AI17").Value = 21
$ws.Range("AJ17").Value = 81
$ws.Range("AM17").Value = 9.5
$ws.Range("AO17").Value = 21
$ws.Range("AR17").Value = 2.05
$ws.Range("AS17").Value = 1.8
$ws.Range("G17").Value = 3.4
$ws.Range("I17").Value = 2.3
$ws.Range("J17").Value = 4.33
$ws.Range("L17").Value = 3.2
$ws.Range("O17").Value = 1.57
$ws.Range("P17").Value = 2.38
$ws.Range("W17").Value = 1.62
$ws.Range("X17").Value = 2.2
$ws.Range("Y17").Value = 2.2
$ws.Range("Z17").Value = 1.62

# Row 18
$ws.Range("AB18").Value = 6.5
$ws.Range("AL18").Value = 13
$ws.Range("G18").Value = 1.57
$ws.Range("H18").Value = 3.8
$ws.Range("I18").Value = 5.75
$ws.Range("J18").Value = 2.2
$ws.Range("N18").Value = 9

# Row 19
$ws.Range("AA19").Value = 7.5
$ws.Range("AC19").Value = 10
$ws.Range("AF19").Value = 34
$ws.Range("AG19").Value = 8.5
$ws.Range("AK19").Value = 351
$ws.Range("AL19").Value = 8
$ws.Range("G19").Value = 2.45
$ws.Range("H19").Value = 3.3
$ws.Range("I19").Value = 2.75
$ws.Range("J19").Value = 3.25
$ws.Range("K19").Value = 2.05
$ws.Range("Y19").Value = 1.83
$ws.Range("Z19").Value = 1.83

# Row 21
$ws.Range("AA21").Value = 7.5
$ws.Range("AC21").Value = 8.5
$ws.Range("AE21").Value = 15
$ws.Range("AF21").Value = 26
$ws.Range("AG21").Value = 10
$ws.Range("AK21").Value = 251
$ws.Range("AL21").Value = 12
$ws.Range("G21").Value = 1.9
$ws.Range("H21").Value = 3.5
$ws.Range("J21").Value = 2.6
$ws.Range("K21").Value = 2.2
$ws.Range("M21").Value = 1.06
$ws.Range("N21").Value = 10
$ws.Range("O21").Value = 1.3
$ws.Range("P21").Value = 3.5
$ws.Range("Q21").Value = 2
$ws.Range("R21").Value = 1.85
$ws.Range("U21").Value = 3.4
$ws.Range("V21").Value = 1.33
$ws.Range("W21").Value = 1.4
$ws.Range("X21").Value = 2.75
$ws.Range("Y21").Value = 1.8
$ws.Range("Z21").Value = 1.91

# Row 22
$ws.Range("AA22").Value = 10
$ws.Range("AB22").Value = 13
$ws.Range("M22").Value = 1.03
$ws.Range("N22").Value = 15
$ws.Range("O22").Value = 1.2
$ws.Range("P22").Value = 4.33
$ws.Range("Q22").Value = 1.65
$ws.Range("R22").Value = 2.2
$ws.Range("U22").Value = 2.63
$ws.Range("V22").Value = 1.44

# Row 25
$ws.Range("AB25").Value = 13
$ws.Range("AC25").Value = 11
$ws.Range("AE25").Value = 23
$ws.Range("AG25").Value = 8
$ws.Range("AK25").Value = 301
$ws.Range("AL25").Value = 8
$ws.Range("AO25").Value = 26
$ws.Range("G25").Value = 2.7
$ws.Range("H25").Value = 3
$ws.Range("I25").Value = 2.7
$ws.Range("J25").Value = 3.4
$ws.Range("K25").Value = 2
$ws.Range("L25").Value = 3.4
$ws.Range("M25").Value = 1.08
$ws.Range("N25").Value = 8
$ws.Range("Q25").Value = 2.2
$ws.Range("R25").Value = 1.65
$ws.Range("U25").Value = 4
$ws.Range("V25").Value = 1.22
$ws.Range("W25").Value = 1.5
$ws.Range("X25").Value = 2.5

# Row 26
$ws.Range("AB26").Value = 9.5
$ws.Range("AD26").Value = 19
$ws.Range("AE26").Value = 19
$ws.Range("AH26").Value = 6
$ws.Range("AL26").Value = 9.5
$ws.Range("AP26").Value = 29
$ws.Range("G26").Value = 2.1
$ws.Range("H26").Value = 3.25
$ws.Range("I26").Value = 3.5
$ws.Range("J26").Value = 2.88
$ws.Range("L26").Value = 4
$ws.Range("M26").Value = 1.07
$ws.Range("N26").Value = 9
